$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.443.89"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "1.968.82"
$ws.Range("E3").Value = "  -5.35%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "244.54"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -4.82%  "
$ws.Range("D7").Value = "58.75"
$ws.Range("E7").Value = "  -9.71%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("D10").Value = "55.80"
$ws.Range("E10").Value = "  -6.51%  "
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "22.01"
$ws.Range("E13").Value = "  -6.64%  "
$ws.Range("D14").Value = "0.833"
$ws.Range("E14").Value = "  -10.60%  "
$ws.Range("D15").Value = "2.255.32"
$ws.Range("E15").Value = "  -5.10%  "
$ws.Range("D16").Value = "13.55"
$ws.Range("E16").Value = "  -8.87%  "
$ws.Range("E17").Value = "  -5.57%  "
$ws.Range("D18").Value = "1.992.66"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("D19").Value = "36.362.73"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").Value = "71.15"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("E22").Value = "  -7.12%  "
$ws.Range("D23").Value = "229.92"
$ws.Range("E23").Value = "  -4.32%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("D28").Value = "164.94"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "19.92"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("D30").Value = "0.123"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("E33").Value = "  -8.36%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").Value = "4.34"
$ws.Range("E35").Value = "  -7.61%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").Value = "5.99"
$ws.Range("E38").Value = "  -7.11%  "
$ws.Range("E39").Value = "  -15.28%  "
$ws.Range("D40").Value = "2.92"
$ws.Range("D41").Value = "0.0967"
$ws.Range("E41").Value = "  -6.42%  "
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("E43").Value = "  -8.54%  "
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("E45").Value = "  -10.04%  "
$ws.Range("D46").Value = "15.81"
$ws.Range("E46").Value = "  -8.33%  "
$ws.Range("D47").Value = "88.96"
$ws.Range("E47").Value = "  -7.31%  "
$ws.Range("D48").Value = "1.351.54"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("E49").Value = "  -8.53%  "
$ws.Range("E50").Value = "  -4.49%  "
$ws.Range("D51").Value = "44.84"
$ws.Range("E51").Value = "  -4.20%  "
